$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Activators sheet: fill in the new "InputBusBAudio" actions for Input 6
# (rows 11-13). Row 11 = Input number, row 12 = Action On, row 13 =
# Action Off.
# ----------------------------------------------------------------------
$wsAct = $wb.Worksheets.Item("Activators")
$wsAct.Range("B11").Value = 1
$wsAct.Range("B12").Value = "red: 0"
$wsAct.Range("B13").Value = "green: 0"

# ----------------------------------------------------------------------
# Shortcuts sheet: button 0 actions become multi-line "string slices"
# and a new button 49 row is added.
# ----------------------------------------------------------------------
$wsShort = $wb.Worksheets.Item("Shortcuts")

$wsShort.Range("B8").Value = "Merge Input=2" + [char]10 + "leds green 0" + [char]10 + "leds yellow 8,9,10"
$wsShort.Range("B8").WrapText = $true

$wsShort.Range("C8").Value = "leds red 0" + [char]10 + "leds yellow 1" + [char]10 + "Merge Input=1"
$wsShort.Range("C8").WrapText = $true

$wsShort.Range("A9").Value = 49
$wsShort.Range("B9").Value = "leds off 49"
$wsShort.Range("B9").Font.Name = "Arial"
$wsShort.Range("B9").Font.Size = 10
$wsShort.Range("B9").NumberFormat = "General"
$wsShort.Range("B9").WrapText = $false
$wsShort.Range("C9").Font.Name = "Arial"
$wsShort.Range("C9").Font.Size = 10
$wsShort.Range("C9").NumberFormat = "General"
$wsShort.Range("C9").WrapText = $false

# ----------------------------------------------------------------------
# Selection / active-sheet bookkeeping so the saved view matches what
# was left selected by the editor.
# ----------------------------------------------------------------------
$wsResp = $wb.Worksheets.Item("Responses")
$wsResp.Range("A17").Select()

$wsAct.Range("A14").Select()

$wsShort.Activate()
$wsShort.Range("B11").Select()
